# feat: add 2022-Q1 data
#
# Before:  Sheets = [ "2020-Q4", "总计" ]
# After:   Sheets = [ "2020-Q4", "2022-Q1", "总计" ]
#
# The existing "总计" sheet (quarter-over-quarter summary) is renamed to
# "2022-Q1" and re-populated with that quarter's per-fund holding detail
# (same shape as the "2020-Q4" sheet). A brand new "总计" sheet is inserted
# right after it, carrying the original summary row plus a new row for
# 2022-Q1 on top.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Re-purpose the current "总计" sheet into the new "2022-Q1" detail
#    sheet. Renaming in place keeps it in slot 2 (right after "2020-Q4"),
#    and its existing header-row / row-label styling is kept for the
#    cells that already carried it.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# E1:H1 are brand-new header cells -> copy the header formatting already
# present on B1 (bold / bordered / centered) onto them.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$q1.Cells.Item(2, 1).Value = 0
$q1.Cells.Item(2, 2).Value = "'002236"
$q1.Cells.Item(2, 3).Value = "大成中证360互联网+大数据100指数A"
$q1.Cells.Item(2, 4).Value = "'5.67"
$q1.Cells.Item(2, 5).Value = "'93.32"
$q1.Cells.Item(2, 6).Value = "'1.05"
$q1.Cells.Item(2, 7).Value = "'0.0595"
$q1.Cells.Item(2, 8).Value = 4

$q1.Cells.Item(3, 1).Value = 1
$q1.Cells.Item(3, 2).Value = "'003359"
$q1.Cells.Item(3, 3).Value = "大成中证360互联网+大数据100指数C"
$q1.Cells.Item(3, 4).Value = "'4.08"
$q1.Cells.Item(3, 5).Value = "'93.32"
$q1.Cells.Item(3, 6).Value = "'1.05"
$q1.Cells.Item(3, 7).Value = "'0.0428"
$q1.Cells.Item(3, 8).Value = 4

# Row 3 is brand new -> copy the row-label styling from A2 onto A3.
$q1.Range("A2").Copy()
$q1.Range("A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Insert a fresh "总计" sheet right after "2022-Q1", holding the
#    quarter-over-quarter summary, newest quarter first.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

# Copy the (already-styled) header formatting from the "2022-Q1" sheet.
$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = 0.1

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2020-Q4"
$total.Cells.Item(3, 3).Value = 2
$total.Cells.Item(3, 4).Value = 0.02

# Copy the row-label styling onto the new column-A cells too.
$q1.Range("A2").Copy()
$total.Range("A2:A3").PasteSpecial(-4122)

# Leave the first sheet ("2020-Q4") selected/active, matching the original.
$wb.Worksheets.Item("2020-Q4").Activate()
